# "convert input to array"
# Sheet1 gains a third column ("الرقم") of numeric data, and a duplicate
# of the sheet ("Sheet1 (2)") is created with the same layout but with the
# names in column A suffixed with "2".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new "الرقم" column to Sheet1 -------------------------------
$ws1.Range("C1").Value = "الرقم"
$ws1.Range("C2").Value = 2000000
$ws1.Range("C3").Value = 1
$ws1.Range("C4").Value = 2
$ws1.Range("C5").Value = 3
$ws1.Range("C6").Value = 4
$ws1.Columns.Item(3).ColumnWidth = 16.71

# --- Duplicate the sheet -> "Sheet1 (2)" --------------------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename the people in column A of the new sheet ---------------------
$ws2.Range("A2").Value = "ابو صلاح2"
$ws2.Range("A3").Value = "لؤي2"
$ws2.Range("A4").Value = "نبيل2"
$ws2.Range("A5").Value = "هشام2"
$ws2.Range("A6").Value = "بورعي2"

# --- Match the recorded selection on both sheets (B2), Sheet1 active ----
$ws2.Range("B2").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B2").Select() | Out-Null
